$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# The "240X120 PORCELANATO" line item was dropped from this advisor's
# breakdown; remove its row and let everything below shift up one row.
$ws.Rows.Item(2).Delete()

# Refresh the remaining line items (labels unchanged except the shift,
# values re-synced from the latest source data) and recompute the
# dependent "POR CUMPLIR" (E = PRESUPUESTO - VENTA) and "CUMPLIMIENTO"
# (F = VENTA / PRESUPUESTO) columns for each row.
$ws.Range("B2").Value = "240X80 PORCELANATO"
$ws.Range("C2").Value = 3120.1145
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 3120.1145
$ws.Range("F2").Value = 0

$ws.Range("B3").Value = "FREGADEROS DE COCINA"
$ws.Range("C3").Value = 646.361575487259
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 646.361575487259
$ws.Range("F3").Value = 0

$ws.Range("B4").Value = "GRANITO"
$ws.Range("C4").Value = 238.32
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 238.32
$ws.Range("F4").Value = 0

$ws.Range("B5").Value = "GRIFERIAS"
$ws.Range("C5").Value = 106.82
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 106.82
$ws.Range("F5").Value = 0

$ws.Range("B6").Value = "INODOROS"
$ws.Range("C6").Value = 1600
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1600
$ws.Range("F6").Value = 0

$ws.Range("B7").Value = "LAVABOS"
$ws.Range("C7").Value = 625
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 625
$ws.Range("F7").Value = 0

$ws.Range("B8").Value = "LED"
$ws.Range("C8").Value = 300
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 300
$ws.Range("F8").Value = 0

$ws.Range("B9").Value = "NO RESURTIBLES"
$ws.Range("C9").Value = 650.25
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 650.25
$ws.Range("F9").Value = 0

$ws.Range("B10").Value = "OTROS"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0

$ws.Range("B11").Value = "PANELES DECORATIVOS"
$ws.Range("C11").Value = 350
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 350
$ws.Range("F11").Value = 0

$ws.Range("B12").Value = "PANELES PU"
$ws.Range("C12").Value = 130
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 130
$ws.Range("F12").Value = 0

$ws.Range("B13").Value = "PANELES PVC"
$ws.Range("C13").Value = 240
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 240
$ws.Range("F13").Value = 0

$ws.Range("B14").Value = "PIEDRA SINTERIZADA"
$ws.Range("C14").Value = 527.03
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 527.03
$ws.Range("F14").Value = 0

$ws.Range("B15").Value = "PORCELANATO"
$ws.Range("C15").Value = 23458.82
$ws.Range("D15").Value = 107.9
$ws.Range("E15").Value = 23350.92
$ws.Range("F15").Value = 0.00459954933794624

$ws.Range("B16").Value = "PUERTAS DE SEGURIDAD"
$ws.Range("C16").Value = 342
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 342
$ws.Range("F16").Value = 0

$ws.Range("B17").Value = "SAL SOLUBLE"
$ws.Range("C17").Value = 1600
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 1600
$ws.Range("F17").Value = 0

# Row 18 is now the TOTAL row; refresh its sums to match the updated
# figures above.
$ws.Range("B18").Value = "TOTAL"
$ws.Range("C18").Value = 33934.71607548726
$ws.Range("D18").Value = 107.9
$ws.Range("E18").Value = 33826.81607548726
$ws.Range("F18").Value = 0.003179634677360438

# Column widths for PRESUPUESTO-adjacent columns D and E were narrowed
# slightly. Excel's ColumnWidth property is expressed in characters and
# is offset from the stored sheet XML width by a fixed padding amount
# (~5/6 of a character), so compensate for that when targeting an exact
# stored width of 11 / 22.
$ws.Columns.Item(4).ColumnWidth = 10.1666666667
$ws.Columns.Item(5).ColumnWidth = 21.1666666667
